# Update NATMI ligand-receptor scores (Rtn4-Cntnap1) with refreshed TPM-based
# expression values. Columns G-J (ligand avg/total expression & specificity),
# K-P (receptor-expressing cells/rate, avg/total expression & specificity) and
# Q-T (edge weights & derived specificities) are recalculated per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 48.46865866666667
$ws.Cells.Item(2, 8).Value = 145.405976
$ws.Cells.Item(2, 9).Value = 0.1554430998624896
$ws.Cells.Item(2, 10).Value = 0.1554430998624896
$ws.Cells.Item(2, 13).Value = 1.432281
$ws.Cells.Item(2, 14).Value = 4.296843
$ws.Cells.Item(2, 15).Value = 0.4004435514722966
$ws.Cells.Item(2, 16).Value = 0.4004435514722965
$ws.Cells.Item(2, 17).Value = 69.420738903752
$ws.Cells.Item(2, 18).Value = 624.786650133768
$ws.Cells.Item(2, 19).Value = 0.06224618696079821
$ws.Cells.Item(2, 20).Value = 0.0622461869607982

$ws.Cells.Item(3, 7).Value = 48.46865866666667
$ws.Cells.Item(3, 8).Value = 145.405976
$ws.Cells.Item(3, 9).Value = 0.1554430998624896
$ws.Cells.Item(3, 10).Value = 0.1554430998624896
$ws.Cells.Item(3, 15).Value = 0.2801347112623808
$ws.Cells.Item(3, 16).Value = 0.2801347112623808
$ws.Cells.Item(3, 17).Value = 48.56404498692267
$ws.Cells.Item(3, 18).Value = 437.0764048823041
$ws.Cells.Item(3, 19).Value = 0.04354500789770796
$ws.Cells.Item(3, 20).Value = 0.04354500789770796

$ws.Cells.Item(4, 7).Value = 48.46865866666667
$ws.Cells.Item(4, 8).Value = 145.405976
$ws.Cells.Item(4, 9).Value = 0.1554430998624896
$ws.Cells.Item(4, 10).Value = 0.1554430998624896
$ws.Cells.Item(4, 13).Value = 0.3795333333333333
$ws.Cells.Item(4, 14).Value = 1.1386
$ws.Cells.Item(4, 15).Value = 0.1061116330539321
$ws.Cells.Item(4, 16).Value = 0.1061116330539321
$ws.Cells.Item(4, 17).Value = 18.39547158595556
$ws.Cells.Item(4, 18).Value = 165.5592442736
$ws.Cells.Item(4, 19).Value = 0.01649432117337423
$ws.Cells.Item(4, 20).Value = 0.01649432117337423

$ws.Cells.Item(5, 7).Value = 48.46865866666667
$ws.Cells.Item(5, 8).Value = 145.405976
$ws.Cells.Item(5, 9).Value = 0.1554430998624896
$ws.Cells.Item(5, 10).Value = 0.1554430998624896
$ws.Cells.Item(5, 13).Value = 0.4135746666666666
$ws.Cells.Item(5, 14).Value = 1.240724
$ws.Cells.Item(5, 15).Value = 0.1156290618384041
$ws.Cells.Item(5, 16).Value = 0.1156290618384041
$ws.Cells.Item(5, 17).Value = 20.04540935184711
$ws.Cells.Item(5, 18).Value = 180.408684166624
$ws.Cells.Item(5, 19).Value = 0.01797373980635304
$ws.Cells.Item(5, 20).Value = 0.01797373980635304

$ws.Cells.Item(6, 7).Value = 48.46865866666667
$ws.Cells.Item(6, 8).Value = 145.405976
$ws.Cells.Item(6, 9).Value = 0.1554430998624896
$ws.Cells.Item(6, 10).Value = 0.1554430998624896
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.3493793333333333
$ws.Cells.Item(6, 14).Value = 1.048138
$ws.Cells.Item(6, 15).Value = 0.0976810423729864
$ws.Cells.Item(6, 16).Value = 0.0976810423729864
$ws.Cells.Item(6, 17).Value = 16.93394765252089
$ws.Cells.Item(6, 18).Value = 152.405528872688
$ws.Cells.Item(6, 19).Value = 0.01518384402425621
$ws.Cells.Item(6, 20).Value = 0.01518384402425621

$ws.Cells.Item(7, 9).Value = 0.20693808715897
$ws.Cells.Item(7, 10).Value = 0.20693808715897
$ws.Cells.Item(7, 13).Value = 1.432281
$ws.Cells.Item(7, 14).Value = 4.296843
$ws.Cells.Item(7, 15).Value = 0.4004435514722966
$ws.Cells.Item(7, 16).Value = 0.4004435514722965
$ws.Cells.Item(7, 17).Value = 92.41835070590598
$ws.Cells.Item(7, 18).Value = 831.765156353154
$ws.Cells.Item(7, 19).Value = 0.08286702255682162
$ws.Cells.Item(7, 20).Value = 0.0828670225568216

$ws.Cells.Item(8, 9).Value = 0.20693808715897
$ws.Cells.Item(8, 10).Value = 0.20693808715897
$ws.Cells.Item(8, 15).Value = 0.2801347112623808
$ws.Cells.Item(8, 16).Value = 0.2801347112623808
$ws.Cells.Item(8, 19).Value = 0.05797054129546746
$ws.Cells.Item(8, 20).Value = 0.05797054129546746

$ws.Cells.Item(9, 9).Value = 0.20693808715897
$ws.Cells.Item(9, 10).Value = 0.20693808715897
$ws.Cells.Item(9, 13).Value = 0.3795333333333333
$ws.Cells.Item(9, 14).Value = 1.1386
$ws.Cells.Item(9, 15).Value = 0.1061116330539321
$ws.Cells.Item(9, 16).Value = 0.1061116330539321
$ws.Cells.Item(9, 17).Value = 24.48949941008889
$ws.Cells.Item(9, 18).Value = 220.4054946908
$ws.Cells.Item(9, 19).Value = 0.02195853836949526
$ws.Cells.Item(9, 20).Value = 0.02195853836949525

$ws.Cells.Item(10, 9).Value = 0.20693808715897
$ws.Cells.Item(10, 10).Value = 0.20693808715897
$ws.Cells.Item(10, 13).Value = 0.4135746666666666
$ws.Cells.Item(10, 14).Value = 1.240724
$ws.Cells.Item(10, 15).Value = 0.1156290618384041
$ws.Cells.Item(10, 16).Value = 0.1156290618384041
$ws.Cells.Item(10, 17).Value = 26.68602640618577
$ws.Cells.Item(10, 18).Value = 240.174237655672
$ws.Cells.Item(10, 19).Value = 0.0239280568768256
$ws.Cells.Item(10, 20).Value = 0.0239280568768256

$ws.Cells.Item(11, 9).Value = 0.20693808715897
$ws.Cells.Item(11, 10).Value = 0.20693808715897
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.3493793333333333
$ws.Cells.Item(11, 14).Value = 1.048138
$ws.Cells.Item(11, 15).Value = 0.0976810423729864
$ws.Cells.Item(11, 16).Value = 0.0976810423729864
$ws.Cells.Item(11, 17).Value = 22.54380373501822
$ws.Cells.Item(11, 18).Value = 202.894233615164
$ws.Cells.Item(11, 19).Value = 0.0202139280603601
$ws.Cells.Item(11, 20).Value = 0.0202139280603601

$ws.Cells.Item(12, 7).Value = 75.47903666666667
$ws.Cells.Item(12, 8).Value = 226.43711
$ws.Cells.Item(12, 9).Value = 0.2420676733554854
$ws.Cells.Item(12, 10).Value = 0.2420676733554854
$ws.Cells.Item(12, 13).Value = 1.432281
$ws.Cells.Item(12, 14).Value = 4.296843
$ws.Cells.Item(12, 15).Value = 0.4004435514722966
$ws.Cells.Item(12, 16).Value = 0.4004435514722965
$ws.Cells.Item(12, 17).Value = 108.10719011597
$ws.Cells.Item(12, 18).Value = 972.9647110437301
$ws.Cells.Item(12, 19).Value = 0.0969344388151064
$ws.Cells.Item(12, 20).Value = 0.09693443881510638

$ws.Cells.Item(13, 7).Value = 75.47903666666667
$ws.Cells.Item(13, 8).Value = 226.43711
$ws.Cells.Item(13, 9).Value = 0.2420676733554854
$ws.Cells.Item(13, 10).Value = 0.2420676733554854
$ws.Cells.Item(13, 15).Value = 0.2801347112623808
$ws.Cells.Item(13, 16).Value = 0.2801347112623808
$ws.Cells.Item(13, 17).Value = 75.62757941082667
$ws.Cells.Item(13, 18).Value = 680.64821469744
$ws.Cells.Item(13, 19).Value = 0.06781155778139522
$ws.Cells.Item(13, 20).Value = 0.06781155778139522

$ws.Cells.Item(14, 7).Value = 75.47903666666667
$ws.Cells.Item(14, 8).Value = 226.43711
$ws.Cells.Item(14, 9).Value = 0.2420676733554854
$ws.Cells.Item(14, 10).Value = 0.2420676733554854
$ws.Cells.Item(14, 13).Value = 0.3795333333333333
$ws.Cells.Item(14, 14).Value = 1.1386
$ws.Cells.Item(14, 15).Value = 0.1061116330539321
$ws.Cells.Item(14, 16).Value = 0.1061116330539321
$ws.Cells.Item(14, 17).Value = 28.64681038288889
$ws.Cells.Item(14, 18).Value = 257.821293446
$ws.Cells.Item(14, 19).Value = 0.02568619612931637
$ws.Cells.Item(14, 20).Value = 0.02568619612931637

$ws.Cells.Item(15, 7).Value = 75.47903666666667
$ws.Cells.Item(15, 8).Value = 226.43711
$ws.Cells.Item(15, 9).Value = 0.2420676733554854
$ws.Cells.Item(15, 10).Value = 0.2420676733554854
$ws.Cells.Item(15, 13).Value = 0.4135746666666666
$ws.Cells.Item(15, 14).Value = 1.240724
$ws.Cells.Item(15, 15).Value = 0.1156290618384041
$ws.Cells.Item(15, 16).Value = 0.1156290618384041
$ws.Cells.Item(15, 17).Value = 31.21621742973778
$ws.Cells.Item(15, 18).Value = 280.94595686764
$ws.Cells.Item(15, 19).Value = 0.02799005797150002
$ws.Cells.Item(15, 20).Value = 0.02799005797150002

$ws.Cells.Item(16, 7).Value = 75.47903666666667
$ws.Cells.Item(16, 8).Value = 226.43711
$ws.Cells.Item(16, 9).Value = 0.2420676733554854
$ws.Cells.Item(16, 10).Value = 0.2420676733554854
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.3493793333333333
$ws.Cells.Item(16, 14).Value = 1.048138
$ws.Cells.Item(16, 15).Value = 0.0976810423729864
$ws.Cells.Item(16, 16).Value = 0.0976810423729864
$ws.Cells.Item(16, 17).Value = 26.37081551124222
$ws.Cells.Item(16, 18).Value = 237.33733960118
$ws.Cells.Item(16, 19).Value = 0.0236454226581674
$ws.Cells.Item(16, 20).Value = 0.0236454226581674

$ws.Cells.Item(17, 7).Value = 51.18999233333333
$ws.Cells.Item(17, 8).Value = 153.569977
$ws.Cells.Item(17, 9).Value = 0.164170647777855
$ws.Cells.Item(17, 10).Value = 0.164170647777855
$ws.Cells.Item(17, 13).Value = 1.432281
$ws.Cells.Item(17, 14).Value = 4.296843
$ws.Cells.Item(17, 15).Value = 0.4004435514722966
$ws.Cells.Item(17, 16).Value = 0.4004435514722965
$ws.Cells.Item(17, 17).Value = 73.31845340917899
$ws.Cells.Item(17, 18).Value = 659.866080682611
$ws.Cells.Item(17, 19).Value = 0.06574107724367174
$ws.Cells.Item(17, 20).Value = 0.06574107724367173

$ws.Cells.Item(18, 7).Value = 51.18999233333333
$ws.Cells.Item(18, 8).Value = 153.569977
$ws.Cells.Item(18, 9).Value = 0.164170647777855
$ws.Cells.Item(18, 10).Value = 0.164170647777855
$ws.Cells.Item(18, 15).Value = 0.2801347112623808
$ws.Cells.Item(18, 16).Value = 0.2801347112623808
$ws.Cells.Item(18, 17).Value = 51.29073423824533
$ws.Cells.Item(18, 18).Value = 461.616608144208
$ws.Cells.Item(18, 19).Value = 0.04598989701300742
$ws.Cells.Item(18, 20).Value = 0.04598989701300742

$ws.Cells.Item(19, 7).Value = 51.18999233333333
$ws.Cells.Item(19, 8).Value = 153.569977
$ws.Cells.Item(19, 9).Value = 0.164170647777855
$ws.Cells.Item(19, 10).Value = 0.164170647777855
$ws.Cells.Item(19, 13).Value = 0.3795333333333333
$ws.Cells.Item(19, 14).Value = 1.1386
$ws.Cells.Item(19, 15).Value = 0.1061116330539321
$ws.Cells.Item(19, 16).Value = 0.1061116330539321
$ws.Cells.Item(19, 17).Value = 19.42830842357777
$ws.Cells.Item(19, 18).Value = 174.8547758122
$ws.Cells.Item(19, 19).Value = 0.01742041553523009
$ws.Cells.Item(19, 20).Value = 0.01742041553523008

$ws.Cells.Item(20, 7).Value = 51.18999233333333
$ws.Cells.Item(20, 8).Value = 153.569977
$ws.Cells.Item(20, 9).Value = 0.164170647777855
$ws.Cells.Item(20, 10).Value = 0.164170647777855
$ws.Cells.Item(20, 13).Value = 0.4135746666666666
$ws.Cells.Item(20, 14).Value = 1.240724
$ws.Cells.Item(20, 15).Value = 0.1156290618384041
$ws.Cells.Item(20, 16).Value = 0.1156290618384041
$ws.Cells.Item(20, 17).Value = 21.17088401592755
$ws.Cells.Item(20, 18).Value = 190.537956143348
$ws.Cells.Item(20, 19).Value = 0.01898289798395645
$ws.Cells.Item(20, 20).Value = 0.01898289798395644

$ws.Cells.Item(21, 7).Value = 51.18999233333333
$ws.Cells.Item(21, 8).Value = 153.569977
$ws.Cells.Item(21, 9).Value = 0.164170647777855
$ws.Cells.Item(21, 10).Value = 0.164170647777855
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 12).Value = 0.6666666666666666
$ws.Cells.Item(21, 13).Value = 0.3493793333333333
$ws.Cells.Item(21, 14).Value = 1.048138
$ws.Cells.Item(21, 15).Value = 0.0976810423729864
$ws.Cells.Item(21, 16).Value = 0.0976810423729864
$ws.Cells.Item(21, 17).Value = 17.88472539475844
$ws.Cells.Item(21, 18).Value = 160.962528552826
$ws.Cells.Item(21, 19).Value = 0.01603636000198928
$ws.Cells.Item(21, 20).Value = 0.01603636000198928

$ws.Cells.Item(22, 7).Value = 72.14667033333333
$ws.Cells.Item(22, 8).Value = 216.440011
$ws.Cells.Item(22, 9).Value = 0.2313804918452
$ws.Cells.Item(22, 10).Value = 0.2313804918452
$ws.Cells.Item(22, 13).Value = 1.432281
$ws.Cells.Item(22, 14).Value = 4.296843
$ws.Cells.Item(22, 15).Value = 0.4004435514722966
$ws.Cells.Item(22, 16).Value = 0.4004435514722965
$ws.Cells.Item(22, 17).Value = 103.334305131697
$ws.Cells.Item(22, 18).Value = 930.008746185273
$ws.Cells.Item(22, 19).Value = 0.09265482589589867
$ws.Cells.Item(22, 20).Value = 0.09265482589589864

$ws.Cells.Item(23, 7).Value = 72.14667033333333
$ws.Cells.Item(23, 8).Value = 216.440011
$ws.Cells.Item(23, 9).Value = 0.2313804918452
$ws.Cells.Item(23, 10).Value = 0.2313804918452
$ws.Cells.Item(23, 15).Value = 0.2801347112623808
$ws.Cells.Item(23, 16).Value = 0.2801347112623808
$ws.Cells.Item(23, 17).Value = 72.28865498054932
$ws.Cells.Item(23, 18).Value = 650.597894824944
$ws.Cells.Item(23, 19).Value = 0.06481770727480278
$ws.Cells.Item(23, 20).Value = 0.06481770727480277

$ws.Cells.Item(24, 7).Value = 72.14667033333333
$ws.Cells.Item(24, 8).Value = 216.440011
$ws.Cells.Item(24, 9).Value = 0.2313804918452
$ws.Cells.Item(24, 10).Value = 0.2313804918452
$ws.Cells.Item(24, 13).Value = 0.3795333333333333
$ws.Cells.Item(24, 14).Value = 1.1386
$ws.Cells.Item(24, 15).Value = 0.1061116330539321
$ws.Cells.Item(24, 16).Value = 0.1061116330539321
$ws.Cells.Item(24, 17).Value = 27.38206628051111
$ws.Cells.Item(24, 18).Value = 246.4385965246
$ws.Cells.Item(24, 19).Value = 0.02455216184651621
$ws.Cells.Item(24, 20).Value = 0.0245521618465162

$ws.Cells.Item(25, 7).Value = 72.14667033333333
$ws.Cells.Item(25, 8).Value = 216.440011
$ws.Cells.Item(25, 9).Value = 0.2313804918452
$ws.Cells.Item(25, 10).Value = 0.2313804918452
$ws.Cells.Item(25, 13).Value = 0.4135746666666666
$ws.Cells.Item(25, 14).Value = 1.240724
$ws.Cells.Item(25, 15).Value = 0.1156290618384041
$ws.Cells.Item(25, 16).Value = 0.1156290618384041
$ws.Cells.Item(25, 17).Value = 29.83803513421822
$ws.Cells.Item(25, 18).Value = 268.542316207964
$ws.Cells.Item(25, 19).Value = 0.02675430919976899
$ws.Cells.Item(25, 20).Value = 0.02675430919976898

$ws.Cells.Item(26, 7).Value = 72.14667033333333
$ws.Cells.Item(26, 8).Value = 216.440011
$ws.Cells.Item(26, 9).Value = 0.2313804918452
$ws.Cells.Item(26, 10).Value = 0.2313804918452
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 12).Value = 0.6666666666666666
$ws.Cells.Item(26, 13).Value = 0.3493793333333333
$ws.Cells.Item(26, 14).Value = 1.048138
$ws.Cells.Item(26, 15).Value = 0.0976810423729864
$ws.Cells.Item(26, 16).Value = 0.0976810423729864
$ws.Cells.Item(26, 17).Value = 25.20655558327978
$ws.Cells.Item(26, 18).Value = 226.859000249518
$ws.Cells.Item(26, 19).Value = 0.02260148762821342
$ws.Cells.Item(26, 20).Value = 0.02260148762821342
